$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - updated daily spot price data
$ws.Range("A2").Value = 46017
$ws.Range("B2").Value = 97.09999999999999
$ws.Range("C2").Value = 91.05
$ws.Range("D2").Value = 85.08
$ws.Range("E2").Value = 80.19
$ws.Range("F2").Value = 72.31
$ws.Range("G2").Value = 73.04000000000001
$ws.Range("H2").Value = 82.72
$ws.Range("I2").Value = 90.86
$ws.Range("J2").Value = 92.94
$ws.Range("K2").Value = 93.19
$ws.Range("L2").Value = 89.25
$ws.Range("M2").Value = 84.98999999999999
$ws.Range("N2").Value = 82.31
$ws.Range("O2").Value = 80.33
$ws.Range("P2").Value = 77.53
$ws.Range("Q2").Value = 78.11
$ws.Range("R2").Value = 85.01000000000001
$ws.Range("S2").Value = 93.45999999999999
$ws.Range("T2").Value = 99.95
$ws.Range("U2").Value = 100.46
$ws.Range("V2").Value = 111.58
$ws.Range("W2").Value = 109.5
$ws.Range("X2").Value = 98.5
$ws.Range("Y2").Value = 91.26000000000001
$ws.Range("Z2").Value = 89.2
$ws.Range("AB2").Value = 102.71
$ws.Range("AD2").Value = 110.54
$ws.Range("AF2").Value = 100.2
$ws.Range("AG2").Value = "2h-16h"
